$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.25"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-3.91%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.89"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-6.44%"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-4.04%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07666"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-6.10%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "4.232"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-2.14%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.613"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-10.89%"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-6.40%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1005"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-9.53%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1723"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-6.73%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08917"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-4.54%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04391"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-5.69%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1055"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.48%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001265"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-1.25%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005816"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.06%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.530"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.53%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.06%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.059"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-5.21%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1342"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-3.53%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "17.57%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.01%"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-4.69%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-5.88%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001221"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-6.42%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.31%"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02345"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "-9.58%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05157"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-5.86%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007919"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.79%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1322"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-5.20%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006548"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-1.13%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.001987"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-6.65%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008120"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-4.35%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3048"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-11.01%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006583"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-5.70%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.26%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003393"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-2.11%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "41.16%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.26%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.26%"
